# Commit: "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# The statement-of-account detail rows (worker doc/name, contribution
# period and overdue salary value) are replaced: the previous two workers'
# records are removed and replaced with updated records for
# "NAVIS JOSE CARO CARO" (9265520) covering periods 2210-2311, followed by
# "JOSE JESUS ORTIZ MARTINEZ" (73145412) covering periods 2206-2311, each
# block sorted most-recent-period-first, with the first row of each block
# carrying a salary value of 37333 instead of 40000.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$rows = @(
    @{ Row = 16; NitDoc = "9265520"; Nombre = "NAVIS JOSE CARO CARO"; Periodo = "2311"; Valor = 37333 },
    @{ Row = 17; NitDoc = "9265520"; Nombre = "NAVIS JOSE CARO CARO"; Periodo = "2310"; Valor = 40000 },
    @{ Row = 18; NitDoc = "9265520"; Nombre = "NAVIS JOSE CARO CARO"; Periodo = "2309"; Valor = 40000 },
    @{ Row = 19; NitDoc = "9265520"; Nombre = "NAVIS JOSE CARO CARO"; Periodo = "2308"; Valor = 40000 },
    @{ Row = 20; NitDoc = "9265520"; Nombre = "NAVIS JOSE CARO CARO"; Periodo = "2307"; Valor = 40000 },
    @{ Row = 21; NitDoc = "9265520"; Nombre = "NAVIS JOSE CARO CARO"; Periodo = "2306"; Valor = 40000 },
    @{ Row = 22; NitDoc = "9265520"; Nombre = "NAVIS JOSE CARO CARO"; Periodo = "2305"; Valor = 40000 },
    @{ Row = 23; NitDoc = "9265520"; Nombre = "NAVIS JOSE CARO CARO"; Periodo = "2304"; Valor = 40000 },
    @{ Row = 24; NitDoc = "9265520"; Nombre = "NAVIS JOSE CARO CARO"; Periodo = "2303"; Valor = 40000 },
    @{ Row = 25; NitDoc = "9265520"; Nombre = "NAVIS JOSE CARO CARO"; Periodo = "2302"; Valor = 40000 },
    @{ Row = 26; NitDoc = "9265520"; Nombre = "NAVIS JOSE CARO CARO"; Periodo = "2301"; Valor = 40000 },
    @{ Row = 27; NitDoc = "9265520"; Nombre = "NAVIS JOSE CARO CARO"; Periodo = "2212"; Valor = 40000 },
    @{ Row = 28; NitDoc = "9265520"; Nombre = "NAVIS JOSE CARO CARO"; Periodo = "2211"; Valor = 40000 },
    @{ Row = 29; NitDoc = "9265520"; Nombre = "NAVIS JOSE CARO CARO"; Periodo = "2210"; Valor = 40000 },
    @{ Row = 30; NitDoc = "73145412"; Nombre = "JOSE JESUS ORTIZ MARTINEZ"; Periodo = "2311"; Valor = 37333 },
    @{ Row = 31; NitDoc = "73145412"; Nombre = "JOSE JESUS ORTIZ MARTINEZ"; Periodo = "2310"; Valor = 40000 },
    @{ Row = 32; NitDoc = "73145412"; Nombre = "JOSE JESUS ORTIZ MARTINEZ"; Periodo = "2309"; Valor = 40000 },
    @{ Row = 33; NitDoc = "73145412"; Nombre = "JOSE JESUS ORTIZ MARTINEZ"; Periodo = "2308"; Valor = 40000 },
    @{ Row = 34; NitDoc = "73145412"; Nombre = "JOSE JESUS ORTIZ MARTINEZ"; Periodo = "2307"; Valor = 40000 },
    @{ Row = 35; NitDoc = "73145412"; Nombre = "JOSE JESUS ORTIZ MARTINEZ"; Periodo = "2306"; Valor = 40000 },
    @{ Row = 36; NitDoc = "73145412"; Nombre = "JOSE JESUS ORTIZ MARTINEZ"; Periodo = "2305"; Valor = 40000 },
    @{ Row = 37; NitDoc = "73145412"; Nombre = "JOSE JESUS ORTIZ MARTINEZ"; Periodo = "2304"; Valor = 40000 },
    @{ Row = 38; NitDoc = "73145412"; Nombre = "JOSE JESUS ORTIZ MARTINEZ"; Periodo = "2303"; Valor = 40000 },
    @{ Row = 39; NitDoc = "73145412"; Nombre = "JOSE JESUS ORTIZ MARTINEZ"; Periodo = "2302"; Valor = 40000 },
    @{ Row = 40; NitDoc = "73145412"; Nombre = "JOSE JESUS ORTIZ MARTINEZ"; Periodo = "2301"; Valor = 40000 },
    @{ Row = 41; NitDoc = "73145412"; Nombre = "JOSE JESUS ORTIZ MARTINEZ"; Periodo = "2212"; Valor = 40000 },
    @{ Row = 42; NitDoc = "73145412"; Nombre = "JOSE JESUS ORTIZ MARTINEZ"; Periodo = "2211"; Valor = 40000 },
    @{ Row = 43; NitDoc = "73145412"; Nombre = "JOSE JESUS ORTIZ MARTINEZ"; Periodo = "2210"; Valor = 40000 },
    @{ Row = 44; NitDoc = "73145412"; Nombre = "JOSE JESUS ORTIZ MARTINEZ"; Periodo = "2209"; Valor = 40000 },
    @{ Row = 45; NitDoc = "73145412"; Nombre = "JOSE JESUS ORTIZ MARTINEZ"; Periodo = "2208"; Valor = 40000 },
    @{ Row = 46; NitDoc = "73145412"; Nombre = "JOSE JESUS ORTIZ MARTINEZ"; Periodo = "2207"; Valor = 40000 },
    @{ Row = 47; NitDoc = "73145412"; Nombre = "JOSE JESUS ORTIZ MARTINEZ"; Periodo = "2206"; Valor = 40000 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.NitDoc
    $ws.Cells.Item($r.Row, 4).Value = $r.Nombre
    $ws.Cells.Item($r.Row, 5).Value = $r.Periodo
    $ws.Cells.Item($r.Row, 6).Value = $r.Valor
}
